$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (order matters for shared-string table layout)
$ws.Range("D1").Value = "Oros"
$ws.Range("E1").Value = "Platas"
$ws.Range("F1").Value = "Bronces"
$ws.Range("A1").Value = "Posicion"

# Move active selection
$ws.Range("B7").Select()
